# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the column headers in row 1 from the generic "_old"/"_new" suffixes
# to the concrete format-version suffixes "_FV2410"/"_FV2504", wraps the data
# range in an Excel Table (ListObject) with an AutoFilter, and freezes the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------------
# Columns A:J previously used the "_old" suffix -> "_FV2410"
# Column K is the neutral "diff" column and is left untouched
# Columns L:U previously used the "_new" suffix -> "_FV2504"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = [string]$cell.Text
    if ($current.EndsWith("_old")) {
        $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2410"
    }
    elseif ($current.EndsWith("_new")) {
        $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2504"
    }
}

# --- 2. Wrap the used range in a Table with AutoFilter ----------------------
$dataRange = $ws.Range("A1:U66")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
